# Applies the commit "Presentation - no 6, 8, 9":
# Collapses several runs that were split around spell-check
# (w:proofErr) marks for "Шорникова" / "разрешени" / "к.ф.-м.н" / "Трунин"
# back into single runs, and fixes a typo
# ("разрешени.." -> "разрешению.") along the way.
#
# Strategy: Find/Replace a span of text that starts at least one
# character before the w:proofErr "spellStart" marker (inside the
# preceding, differently-formatted run) and ends at least one character
# after the matching "spellEnd" marker. Replacing that whole span
# collapses the runs inside it into one run and drops the now-orphaned
# w:proofErr markers.

$d = $word.ActiveDocument
$rng = $d.Content

# --- Location 1: "Шорникова Александра Евгеньевича," -----------------
$old1 = "Института математики и информатики, `vШорникова Александра Евгеньевича,"
$new1 = "Института математики и информатики, `vШорникова Александра Евгеньевича,"
$r1 = $rng.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)
Write-Output "loc1: $r1"

# --- Location 2: "А. Е. Шорникова посвящена ... бросания лучей." -----
$old2 = "Выпускная квалификационная работа А. Е. Шорникова посвящена рассмотрению проблем отсутствия адекватных способов создания интерактивных планов помещений и предлагается свой выход из сложившейся ситуации. Автор создал свою программную реализацию такого интерактивного плана используя метод бросания лучей."
$new2 = "Выпускная квалификационная работа А. Е. Шорникова посвящена рассмотрению проблем отсутствия адекватных способов создания интерактивных планов помещений и предлагается свой выход из сложившейся ситуации. Автор создал свою программную реализацию такого интерактивного плана используя метод бросания лучей."
$r2 = $rng.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)
Write-Output "loc2: $r2"

# --- Location 3: "...подлежащих разрешени.." -> "...разрешению." -----
$old3 = "сферы IT-рынка, очерчивает круг проблем, подлежащих разрешени.. В основных частях работы последовательно рассматриваются вопросы, обосновывается выбор тех или иных средств и методов, последовательно рассматриваются важные для рассмотрения темы и выводятся результаты работы."
$new3 = "сферы IT-рынка, очерчивает круг проблем, подлежащих разрешению. В основных частях работы последовательно рассматриваются вопросы, обосновывается выбор тех или иных средств и методов, последовательно рассматриваются важные для рассмотрения темы и выводятся результаты работы."
$r3 = $rng.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2)
Write-Output "loc3: $r3"

# --- Location 4: "Представленная на защиту работа А. Е. Шорникова соответствует..." ---
$old4 = "Представленная на защиту работа А. Е. Шорникова соответствует предъявляемым требованиям, а ее автор заслуживает положительной оценки."
$new4 = "Представленная на защиту работа А. Е. Шорникова соответствует предъявляемым требованиям, а ее автор заслуживает положительной оценки."
$r4 = $rng.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2)
Write-Output "loc4: $r4"

# --- Location 5: "Научный руководитель: к.ф.-м.н., ст. преп.  " ------
$old5 = "Научный руководитель: к.ф.-м.н., ст. преп.  "
$new5 = "Научный руководитель: к.ф.-м.н., ст. преп.  "
$r5 = $rng.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, $new5, 2)
Write-Output "loc5: $r5"

# --- Location 6: "Д. О. Трунин" ---------------------------------------
$old6 = "Д. О. Трунин"
$new6 = "Д. О. Трунин"
$r6 = $rng.Find.Execute($old6, $true, $false, $false, $false, $false, $true, 1, $false, $new6, 2)
Write-Output "loc6: $r6"
